$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing GDP values (column B) for rows 12-30 (revised figures)
$ws.Range("B12").Value = 80171.426000000007
$ws.Range("B13").Value = 85318.547000000006
$ws.Range("B14").Value = 91144.945000000007
$ws.Range("B15").Value = 99483.963000000003
$ws.Range("B16").Value = 107620.219
$ws.Range("B17").Value = 113539.03200000001
$ws.Range("B18").Value = 117748.606
$ws.Range("B19").Value = 115724.10400000001
$ws.Range("B20").Value = 113934.961
$ws.Range("B21").Value = 115728.996
$ws.Range("B22").Value = 117019.432
$ws.Range("B23").Value = 121305.114
$ws.Range("B24").Value = 126082.901
$ws.Range("B25").Value = 131043.519
$ws.Range("B26").Value = 139192.47200000001
$ws.Range("B27").Value = 145593.76800000001
$ws.Range("B28").Value = 150980.261
$ws.Range("B29").Value = 159127.777
$ws.Range("B30").Value = 168393.55799999999

# Append the new 2020-01-01 observation as row 31, copying row 30's
# formatting (date format in A, number format in B) onto it first.
$ws.Range("A30:B30").Copy()
$ws.Range("A31").PasteSpecial(-4122)

$ws.Range("A31").Value = 43831
$ws.Range("B31").Value = 169269.54

# Reflect the new full-column selection recorded on the sheet view.
$ws.Range("A1:B1048576").Select()
